$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the date on row 7 (was 21-Feb-2023, now 21-Jan-2023)
$ws.Range("E7").Value = (Get-Date -Year 2023 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0).Date

# Add closing quantity row 8 (set number formats before values so the
# existing style entries get reused instead of new ones being created)
$ws.Range("A8").Value = "AD480"

$ws.Range("B8").Value = 1

$ws.Range("C8").Value = 55

$ws.Range("D8").NumberFormat = "0.00"
$ws.Range("D8").Value = 35

$ws.Range("E8").NumberFormat = "d-mmm-yy"
$ws.Range("E8").Value = (Get-Date -Year 2023 -Month 1 -Day 22 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("E8").Select()
